$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @(0.244, 0.093, 0.046, -48, -1, -153),
    @(0.227, -0.172, -0.357, -87, 6, -42),
    @(0.207, 0.437, -0.08599999999999999, -173, 14, 51),
    @(0.189, -0.325, -0.047, 103, 35, -13),
    @(0.037, -0.187, 0.385, 27, -28, 98),
    @(-0.064, -0.11, 0.492, 35, 18, 0),
    @(-0.134, 0.043, 0.359, -30, 15, 93),
    @(0.453, 0.122, 0.185, 44, 23, 3),
    @(-0.276, 0.079, 0.101, -7, 2, 28),
    @(0.2, 0.035, -0.16, -13, 54, -62)
)

for ($r = 1; $r -le 10; $r++) {
    for ($c = 1; $c -le 6; $c++) {
        $ws.Cells.Item($r, $c).Value = $values[$r - 1][$c - 1]
    }
}
